# Reserva sheet: add a "Name" column after RES_Id, and append mocked
# vital-signs / triage columns (Age, Blood, temperature, heartbeat, wound,
# manchester, body_part) for the 7 patient rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reserva")

# --- Insert the new "Name" column as column B (shifts old B..G to C..H) ---
$ws.Columns("B:B").Insert()

$ws.Range("B1").Value = "Name"
$ws.Range("B2").Value = "Patient 1"
$ws.Range("B3").Value = "Patient 2"
$ws.Range("B4").Value = "Patient 3"
$ws.Range("B5").Value = "Patient 4"
$ws.Range("B6").Value = "Patient 5"
$ws.Range("B7").Value = "Patient 6"
$ws.Range("B8").Value = "Patient 7"

# --- Append new data columns I..O ---
$ws.Range("I1").Value = "Age"
$ws.Range("J1").Value = "Blood"
$ws.Range("K1").Value = "temperature"
$ws.Range("L1").Value = "heartbeat"
$ws.Range("M1").Value = "wound"
$ws.Range("N1").Value = "manchester"
$ws.Range("O1").Value = "body_part"

$ws.Range("I2").Value = 30
$ws.Range("J2").Value = "O+"
$ws.Range("K2").Value = 36.6
$ws.Range("L2").Value = 80
$ws.Range("M2").Value = "Escoriação"
$ws.Range("N2").Value = 2
$ws.Range("O2").Value = "Membros Inferiores"

$ws.Range("I3").Value = 40
$ws.Range("J3").Value = "A+"
$ws.Range("K3").Value = 36.5
$ws.Range("L3").Value = 75
$ws.Range("M3").Value = "Fratura"
$ws.Range("N3").Value = 3
$ws.Range("O3").Value = "Membros Superiores"

$ws.Range("I4").Value = 50
$ws.Range("J4").Value = "B-"
$ws.Range("K4").Value = 38.5
$ws.Range("L4").Value = 75
$ws.Range("M4").Value = "Contusão"
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = "Coluna"

$ws.Range("I5").Value = 38
$ws.Range("J5").Value = "AB+"
$ws.Range("K5").Value = 36.6
$ws.Range("L5").Value = 90
$ws.Range("M5").Value = "Traumatismo"
$ws.Range("N5").Value = 4
$ws.Range("O5").Value = "Cabeça"

$ws.Range("I6").Value = 42
$ws.Range("J6").Value = "O-"
$ws.Range("K6").Value = 36.5
$ws.Range("L6").Value = 88
$ws.Range("M6").Value = "Hematoma"
$ws.Range("N6").Value = 3
$ws.Range("O6").Value = "Abdômen"

$ws.Range("I7").Value = 25
$ws.Range("J7").Value = "O+"
$ws.Range("K7").Value = 40
$ws.Range("L7").Value = 80
$ws.Range("M7").Value = "Luxação"
$ws.Range("N7").Value = 2
$ws.Range("O7").Value = "Pelve"

$ws.Range("I8").Value = 66
$ws.Range("J8").Value = "B+"
$ws.Range("K8").Value = 36.6
$ws.Range("L8").Value = 120
$ws.Range("M8").Value = "Traumatismo"
$ws.Range("N8").Value = 5
$ws.Range("O8").Value = "Pescoço"

# --- Formatting: remove the old per-row wrap-height, widen the new
#     date columns, bestFit the manchester column ---
$ws.Rows("2:8").RowHeight = 14.4

$ws.Columns("E:G").ColumnWidth = 13.44140625
$ws.Columns("N:N").ColumnWidth = 10.5546875

# --- Autofilter over the now-wider table, with the two active filters on
#     Name (col 2) and ESP_ID (col 7) that the mock data exercises ---
$ws.Range("A1:H8").AutoFilter(2)
$ws.Range("A1:H8").AutoFilter(7)

# --- The PivotTable/filter-tracking defined name grows with the table ---
$wb.Names.Add("_xlnm._FilterDatabase", "=Reserva!$A$1:$H$8")
